$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 (shifts existing rows 11-27 down to 12-28,
# copying formatting from the row above as Excel does by default).
$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the new weekly price observation.
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Vega Modelo de Temuco"
$ws.Range("C11").Value = "La Araucanía"
$ws.Range("D11").Value = 45225
$ws.Range("E11").Value = 9
$ws.Range("F11").Value = 100112036
$ws.Range("G11").Value = "Caigua"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 80
$ws.Range("K11").Value = 20000
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = 20000
$ws.Range("N11").Value = "$/caja 15 kilos"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 1333
$ws.Range("Q11").Value = 15
$ws.Range("R11").Value = "Hortaliza"
